$d = $word.ActiveDocument

# --- Locate the block to remove -------------------------------------------------
# Sequence right after the (red) stack-trace run is:
#   [run: "    "] [run: "<---" (orange)] [run: "M2Doc version mismatch..." (orange)] [run: "    "] [run: "demonstration"]
# We need to delete the first three runs in that list and leave the trailing
# "    " + "demonstration" runs exactly as they were (two separate, unformatted
# runs).

# Unique anchor: the last line of the (single) stack trace in the document.
$afterStack = $d.Content
$null = $afterStack.Find.Execute("Thread.run(Thread.java:748)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$runBoundary = $afterStack.End + 1   # the run's text ends with a trailing newline

# The second "<---" marker (the first one belongs to the earlier red block).
$arrowRng = $d.Range($runBoundary, $d.Content.End)
$null = $arrowRng.Find.Execute("<---", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# The "M2Doc version mismatch..." run that follows the arrow.
$msgRng = $d.Range($arrowRng.End, $d.Content.End)
$null = $msgRng.Find.Execute("M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# The leading 4-space run right after the stack trace, before the arrow.
$spacesRng = $d.Range($runBoundary, $arrowRng.Start)

# --- Remove the three runs (rightmost first so earlier offsets stay valid) -----
$msgRng.Delete()
$arrowRng.Delete()
$spacesRng.Delete()

# --- Repair the trailing "    " / "demonstration" split ------------------------
# Deleting text anywhere in this paragraph makes the engine re-coalesce any
# directly-adjacent, identically-formatted runs - which merges the untouched
# "    " and "demonstration" runs into one. Re-split them back into two plain
# runs (matching the original, untouched document) by re-inserting
# "demonstration" as its own run via InsertXML (this does not leave a stray
# empty <w:rPr/> the way toggling a character-formatting property would).
$demoRng = $d.Content
$null = $demoRng.Find.Execute("demonstration", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$demoXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>demonstration</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$demoRng.InsertXML($demoXml)
